$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 39, which shifts existing rows 39-42 down to 40-43
$ws.Rows("39").Insert()

# Populate the new row 39 with the new data record
$ws.Cells.Item(39, 1).Value = 11
$ws.Cells.Item(39, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(39, 3).Value = "Bíobío"
$ws.Cells.Item(39, 4).Value = 44461
$ws.Cells.Item(39, 5).Value = 8
$ws.Cells.Item(39, 6).Value = 100112012
$ws.Cells.Item(39, 7).Value = "Espinaca"
$ws.Cells.Item(39, 8).Value = "Sin especificar"
$ws.Cells.Item(39, 9).Value = "Primera"
$ws.Cells.Item(39, 10).Value = 50
$ws.Cells.Item(39, 11).Value = 7000
$ws.Cells.Item(39, 12).Value = 7500
$ws.Cells.Item(39, 13).Value = 7300
$ws.Cells.Item(39, 14).Value = "$/cuna 10 kilos"
$ws.Cells.Item(39, 15).Value = "Provincia de Chacabuco"
$ws.Cells.Item(39, 16).Value = 730
$ws.Cells.Item(39, 17).Value = 10
$ws.Cells.Item(39, 18).Value = "Hortaliza"
